# Refresh crypto snapshot: update Price (column D) and Volume(1h) (column E)
# text values to match the latest scrape, as produced by the GitHub Actions job.
#
# D-column price strings look numeric (e.g. "291.64") but must stay plain TEXT
# (matching the original inlineStr cells) - Excel auto-coerces a bare numeric-
# looking .Value into a real number, so we force the cell to Text format first,
# assign the literal string, then restore the default "Normal" style so no
# visible formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.486.74"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.39"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.64"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3701"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  -2.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.91"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3386"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07549"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.26"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.031"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.961"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.582.75"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +1.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001121"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.74"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06758"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.300"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.43"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -1.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.21"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.499.10"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.370"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -1.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.602"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  -3.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.06"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.14"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.055"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.19"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.755.60"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("E32").Value = "  +7.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.253"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.017"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.782"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08355"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  -1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02487"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("E38").Value = "  -4.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2300"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06538"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.446"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.35"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6223"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  -1.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.05"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.807"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5859"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.24"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +3.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.071"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07331"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +0.11%  "
